$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume/coin data per latest scrape
$ws.Range("D2").Value = "'" + '57.937.01'
$ws.Range("E2").Value = '  +2.54%  '
$ws.Range("D3").Value = "'" + '3.048.75'
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'" + '519.61'
$ws.Range("E5").Value = '  +3.60%  '
$ws.Range("D6").Value = "'" + '141.17'
$ws.Range("E6").Value = '  +4.68%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = "'" + '0.446'
$ws.Range("E8").Value = '  +4.27%  '
$ws.Range("D9").Value = "'" + '7.50'
$ws.Range("E9").Value = '  +2.95%  '
$ws.Range("D10").Value = "'" + '0.111'
$ws.Range("E10").Value = '  +4.17%  '
$ws.Range("D11").Value = "'" + '0.367'
$ws.Range("E11").Value = '  +4.40%  '
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = "'" + '3.571.32'
$ws.Range("E12").Value = '  +2.35%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = "'" + '0.131'
$ws.Range("E13").Value = '  +2.98%  '
$ws.Range("D14").Value = "'" + '26.89'
$ws.Range("E14").Value = '  +7.25%  '
$ws.Range("D15").Value = "'" + '0.0000170'
$ws.Range("E15").Value = '  +11.10%  '
$ws.Range("D16").Value = "'" + '57.906.68'
$ws.Range("E16").Value = '  +2.61%  '
$ws.Range("D17").Value = "'" + '6.22'
$ws.Range("E17").Value = '  +9.40%  '
$ws.Range("D18").Value = "'" + '3.049.05'
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("D19").Value = "'" + '13.10'
$ws.Range("E19").Value = '  +6.40%  '
$ws.Range("D20").Value = "'" + '8.12'
$ws.Range("E20").Value = '  +4.24%  '
$ws.Range("D21").Value = "'" + '338.27'
$ws.Range("E21").Value = '  +2.91%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = "'" + '5.72'
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("D24").Value = "'" + '0.501'
$ws.Range("E24").Value = '  +6.57%  '
$ws.Range("D25").Value = "'" + '65.15'
$ws.Range("E25").Value = '  +4.52%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = "'" + '0.167'
$ws.Range("E26").Value = '  +2.53%  '
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = "'" + '0.0₃0955'
$ws.Range("E27").Value = '  +5.97%  '
$ws.Range("D28").Value = "'" + '1.01'
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("D29").Value = "'" + '6.88'
$ws.Range("E29").Value = '  +6.64%  '
$ws.Range("D30").Value = "'" + '7.41'
$ws.Range("E30").Value = '  +8.17%  '
$ws.Range("D31").Value = "'" + '1.83'
$ws.Range("E31").Value = '  +4.67%  '
$ws.Range("D32").Value = "'" + '1.21'
$ws.Range("E32").Value = '  +4.05%  '
$ws.Range("D33").Value = "'" + '21.06'
$ws.Range("E33").Value = '  +3.31%  '
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").Value = "'" + '158.98'
$ws.Range("E34").Value = '  +2.61%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Value = "'" + '4.74'
$ws.Range("E35").Value = '  +6.56%  '
$ws.Range("D36").Value = "'" + '5.91'
$ws.Range("E36").Value = '  +6.05%  '
$ws.Range("D37").Value = "'" + '1.30'
$ws.Range("E37").Value = '  +0.80%  '
$ws.Range("D38").Value = "'" + '25.43'
$ws.Range("E38").Value = '  +10.87%  '
$ws.Range("D39").Value = "'" + '0.0694'
$ws.Range("E39").Value = '  +3.05%  '
$ws.Range("D40").Value = "'" + '3.083.92'
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("D41").Value = "'" + '37.75'
$ws.Range("E41").Value = '  +3.78%  '
$ws.Range("D42").Value = "'" + '3.92'
$ws.Range("E42").Value = '  +9.75%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = "'" + '1.00'
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = "'" + '0.664'
$ws.Range("E44").Value = '  +3.60%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = "'" + '2.328.98'
$ws.Range("E45").Value = '  +4.25%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = "'" + '1.45'
$ws.Range("E46").Value = '  +4.11%  '
$ws.Range("D47").Value = "'" + '1.02'
$ws.Range("E47").Value = '  +2.15%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = "'" + '6.04'
$ws.Range("E48").Value = '  +4.83%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = "'" + '0.0241'
$ws.Range("E49").Value = '  +2.42%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = "'" + '19.80'
$ws.Range("E50").Value = '  +4.42%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = "'" + '1.91'
$ws.Range("E51").Value = '  -2.31%  '
